$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 176; this shifts rows 176..217 down to 177..218
# and keeps the date-column (D) number format on the new row.
$ws.Rows("176:176").Insert()

# Populate the newly inserted row 176 with the new weekly record.
$ws.Cells.Item(176, 1).Value = 5
$ws.Cells.Item(176, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(176, 3).Value = "Maule"
$ws.Cells.Item(176, 4).Value = 44785
$ws.Cells.Item(176, 5).Value = 7
$ws.Cells.Item(176, 6).Value = 100112017
$ws.Cells.Item(176, 7).Value = "Apio"
$ws.Cells.Item(176, 8).Value = "Americana (o)"
$ws.Cells.Item(176, 9).Value = "Primera"
$ws.Cells.Item(176, 10).Value = 800
$ws.Cells.Item(176, 11).Value = 10000
$ws.Cells.Item(176, 12).Value = 10000
$ws.Cells.Item(176, 13).Value = 10000
$ws.Cells.Item(176, 14).Value = "`$/docena de matas"
$ws.Cells.Item(176, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(176, 16).Value = 1667
$ws.Cells.Item(176, 17).Value = 6
$ws.Cells.Item(176, 18).Value = "Hortaliza"
